# Auto-generated edit script: apply value changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.4
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 2.2
$ws.Range("O2").Value = 1.67
$ws.Range("AA2").Value = 9
# Row 3
$ws.Range("G3").Value = 2.5
$ws.Range("H3").Value = 3.2
$ws.Range("N3").Value = 2.2
$ws.Range("O3").Value = 1.67
$ws.Range("Y3").Value = 41
$ws.Range("AE3").Value = 9
# Row 4
$ws.Range("G4").Value = 2.3
$ws.Range("I4").Value = 3.1
$ws.Range("N4").Value = 2.02
$ws.Range("O4").Value = 1.88
$ws.Range("AG4").Value = 12
# Row 5
$ws.Range("G5").Value = 1.29
$ws.Range("H5").Value = 6
$ws.Range("L5").Value = 1.13
$ws.Range("M5").Value = 6
$ws.Range("W5").Value = 9
$ws.Range("AA5").Value = 12
$ws.Range("AE5").Value = 23
# Row 6
$ws.Range("N6").Value = 1.75
$ws.Range("O6").Value = 2.05
# Row 7
$ws.Range("G7").Value = 1.73
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 4.33
$ws.Range("T7").Value = 9.5
$ws.Range("U7").Value = 9.5
$ws.Range("W7").Value = 15
$ws.Range("AA7").Value = 7.5
$ws.Range("AD7").Value = 126
$ws.Range("AF7").Value = 26
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 41
# Row 8
$ws.Range("G8").Value = 2.63
$ws.Range("I8").Value = 2.9
$ws.Range("N8").Value = 2.7
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 1.62
$ws.Range("Q8").Value = 2.2
$ws.Range("V8").Value = 11
$ws.Range("W8").Value = 26
# Row 9
$ws.Range("G9").Value = 3.5
$ws.Range("I9").Value = 2.2
$ws.Range("Y9").Value = 51
$ws.Range("AE9").Value = 5.5
$ws.Range("AI9").Value = 21
# Row 10
$ws.Range("G10").Value = 1.9
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 3.9
$ws.Range("L10").Value = 1.25
$ws.Range("M10").Value = 3.75
$ws.Range("N10").Value = 1.83
$ws.Range("O10").Value = 2.03
$ws.Range("W10").Value = 17
$ws.Range("AA10").Value = 6.5
$ws.Range("AD10").Value = 201
# Row 12
$ws.Range("J12").Value = 1.06
$ws.Range("K12").Value = 10
$ws.Range("L12").Value = 1.33
$ws.Range("M12").Value = 3.25
# Row 13
$ws.Range("I13").Value = 1.67
$ws.Range("J13").Value = 1.07
$ws.Range("K13").Value = 9
$ws.Range("R13").Value = 2.2
$ws.Range("S13").Value = 1.62
$ws.Range("AB13").Value = 21
$ws.Range("AC13").Value = 81
$ws.Range("AH13").Value = 12
# Row 14
$ws.Range("H14").Value = 3.6
$ws.Range("K14").Value = 17
$ws.Range("N14").Value = 1.57
$ws.Range("O14").Value = 2.35
$ws.Range("P14").Value = 1.29
$ws.Range("Q14").Value = 3.5
$ws.Range("U14").Value = 19
$ws.Range("Z14").Value = 17
$ws.Range("AG14").Value = 9.5
# Row 15
$ws.Range("K15").Value = 9
# Row 16
$ws.Range("G16").Value = 3.2
$ws.Range("I16").Value = 2.2
$ws.Range("J16").Value = 1.07
$ws.Range("K16").Value = 9
$ws.Range("N16").Value = 2.1
$ws.Range("O16").Value = 1.7
$ws.Range("R16").Value = 1.91
$ws.Range("S16").Value = 1.91
$ws.Range("T16").Value = 9
$ws.Range("U16").Value = 15
$ws.Range("Y16").Value = 41
$ws.Range("AE16").Value = 7
$ws.Range("AG16").Value = 9.5
# Row 18
$ws.Range("G18").Value = 1.83
$ws.Range("H18").Value = 3.4
$ws.Range("I18").Value = 4.33
$ws.Range("T18").Value = 5.5
$ws.Range("U18").Value = 7.5
$ws.Range("X18").Value = 17
$ws.Range("AB18").Value = 21
$ws.Range("AC18").Value = 81
$ws.Range("AE18").Value = 9.5
$ws.Range("AF18").Value = 21
$ws.Range("AH18").Value = 51
$ws.Range("AJ18").Value = 51
# Row 20
$ws.Range("J20").Value = 1.06
$ws.Range("K20").Value = 10
$ws.Range("R20").Value = 2.38
$ws.Range("S20").Value = 1.53
$ws.Range("Z20").Value = 10
$ws.Range("AC20").Value = 101
$ws.Range("AE20").Value = 17
# Row 21
$ws.Range("G21").Value = 7
$ws.Range("I21").Value = 1.38
# Row 22
$ws.Range("H22").Value = 3.5
$ws.Range("I22").Value = 2.3
$ws.Range("J22").Value = 1.05
$ws.Range("K22").Value = 11
$ws.Range("L22").Value = 1.3
$ws.Range("M22").Value = 3.4
$ws.Range("N22").Value = 2
$ws.Range("O22").Value = 1.85
$ws.Range("P22").Value = 1.44
$ws.Range("Q22").Value = 2.63
$ws.Range("R22").Value = 1.8
$ws.Range("S22").Value = 1.91
$ws.Range("T22").Value = 9
$ws.Range("Y22").Value = 34
$ws.Range("Z22").Value = 10
$ws.Range("AA22").Value = 6.5
$ws.Range("AC22").Value = 51
$ws.Range("AD22").Value = 251
$ws.Range("AE22").Value = 7.5
$ws.Range("AI22").Value = 19
$ws.Range("AJ22").Value = 29
# Row 23
$ws.Range("G23").Value = 2.8
$ws.Range("I23").Value = 2.3
$ws.Range("AD23").Value = 201
$ws.Range("AH23").Value = 21
# Row 24
$ws.Range("G24").Value = 3.75
$ws.Range("H24").Value = 3.3
$ws.Range("I24").Value = 2.1
$ws.Range("J24").Value = 1.08
$ws.Range("K24").Value = 8
$ws.Range("L24").Value = 1.44
$ws.Range("M24").Value = 2.63
$ws.Range("N24").Value = 2.35
$ws.Range("O24").Value = 1.57
$ws.Range("P24").Value = 1.53
$ws.Range("Q24").Value = 2.38
$ws.Range("R24").Value = 2.05
$ws.Range("S24").Value = 1.7
$ws.Range("T24").Value = 9
$ws.Range("U24").Value = 17
$ws.Range("Z24").Value = 7.5
$ws.Range("AB24").Value = 17
$ws.Range("AD24").Value = 451
$ws.Range("AF24").Value = 9
$ws.Range("AH24").Value = 19
$ws.Range("AI24").Value = 19
$ws.Range("AJ24").Value = 34
# Row 25
$ws.Range("G25").Value = 1.75
$ws.Range("H25").Value = 3.7
$ws.Range("I25").Value = 4.5
$ws.Range("N25").Value = 1.93
$ws.Range("O25").Value = 1.93
$ws.Range("T25").Value = 7
$ws.Range("U25").Value = 8.5
$ws.Range("W25").Value = 15
$ws.Range("Z25").Value = 11
$ws.Range("AA25").Value = 7
$ws.Range("AE25").Value = 12
$ws.Range("AF25").Value = 23
$ws.Range("AG25").Value = 15
$ws.Range("AH25").Value = 51
# Row 26
$ws.Range("G26").Value = 2.75
$ws.Range("I26").Value = 2.2
$ws.Range("R26").Value = 1.75
$ws.Range("S26").Value = 2
$ws.Range("U26").Value = 15
$ws.Range("W26").Value = 29
$ws.Range("X26").Value = 23
$ws.Range("Y26").Value = 29
$ws.Range("AE26").Value = 8
$ws.Range("AF26").Value = 11
$ws.Range("AG26").Value = 9.5
$ws.Range("AH26").Value = 21
# Row 27
$ws.Range("G27").Value = 2.6
$ws.Range("I27").Value = 2.55
$ws.Range("U27").Value = 15
$ws.Range("V27").Value = 10
$ws.Range("W27").Value = 26
$ws.Range("X27").Value = 19
$ws.Range("Z27").Value = 13
$ws.Range("AE27").Value = 11
$ws.Range("AG27").Value = 10
$ws.Range("AH27").Value = 26
$ws.Range("AI27").Value = 19
# Row 28
$ws.Range("N28").Value = 1.9
$ws.Range("O28").Value = 1.95
# Row 29
$ws.Range("G29").Value = 1.98
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 3.45
$ws.Range("K29").Value = 7.9
$ws.Range("L29").Value = 1.24
$ws.Range("M29").Value = 3.6
$ws.Range("N29").Value = 1.72
$ws.Range("O29").Value = 2
$ws.Range("P29").Value = 1.36
$ws.Range("Q29").Value = 2.87
$ws.Range("R29").Value = 1.65
$ws.Range("S29").Value = 2.12
$ws.Range("T29").Value = 8.5
$ws.Range("U29").Value = 10.25
$ws.Range("Y29").Value = 23
$ws.Range("Z29").Value = 7.9
$ws.Range("AA29").Value = 6.8
$ws.Range("AB29").Value = 13
$ws.Range("AD29").Value = 350
$ws.Range("AE29").Value = 11.75
$ws.Range("AF29").Value = 19.5
$ws.Range("AJ29").Value = 32
